$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.670.20"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "2.647.85"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'595.44"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("D6").Value = "'156.33"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'0.625"
$ws.Range("E8").Value = "  +1.29%  "
$ws.Range("D9").Value = "'0.126"
$ws.Range("E9").Value = "  +2.31%  "
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").Value = "'5.78"
$ws.Range("E10").Value = "  -1.61%  "
$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").Value = "'0.393"
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("E12").Value = "  +1.36%  "
$ws.Range("D13").Value = "'0.0000196"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").Value = "'28.35"
$ws.Range("E14").Value = "  -3.06%  "
$ws.Range("D15").Value = "3.123.68"
$ws.Range("E15").Value = "  -0.49%  "
$ws.Range("D16").Value = "65.501.00"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "2.648.38"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").Value = "'12.50"
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("D19").Value = "'4.70"
$ws.Range("E19").Value = "  -1.85%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'347.56"
$ws.Range("E20").Value = "  -0.84%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'7.40"
$ws.Range("E21").Value = "  -1.95%  "
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").Value = "'68.87"
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("D24").Value = "'0.0000112"
$ws.Range("E24").Value = "  +2.61%  "
$ws.Range("E25").Value = "  +4.86%  "
$ws.Range("D26").Value = "'9.53"
$ws.Range("E26").Value = "  -1.63%  "
$ws.Range("D27").Value = "'1.57"
$ws.Range("E27").Value = "  -1.01%  "
$ws.Range("D28").Value = "'0.163"
$ws.Range("E28").Value = "  -2.43%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").Value = "'7.87"
$ws.Range("E30").Value = "  -2.42%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").Value = "'533.49"
$ws.Range("E31").Value = "  -1.61%  "
$ws.Range("D32").Value = "'2.12"
$ws.Range("E32").Value = "  -0.71%  "
$ws.Range("D33").Value = "'1.74"
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("D34").Value = "'6.37"
$ws.Range("E34").Value = "  -1.89%  "
$ws.Range("D35").Value = "'5.38"
$ws.Range("E35").Value = "  -0.46%  "
$ws.Range("D36").Value = "'0.416"
$ws.Range("E36").Value = "  -1.17%  "
$ws.Range("D37").Value = "'20.22"
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("D38").Value = "'0.998"
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "'155.34"
$ws.Range("E39").Value = "  -2.55%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'1.92"
$ws.Range("E40").Value = "  -0.90%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").Value = "'161.19"
$ws.Range("E42").Value = "  -2.09%  "
$ws.Range("D43").Value = "'4.03"
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("D44").Value = "'2.27"
$ws.Range("E44").Value = "  -0.76%  "
$ws.Range("D45").Value = "'0.0599"
$ws.Range("E45").Value = "  -1.26%  "
$ws.Range("D46").Value = "'22.34"
$ws.Range("E46").Value = "  -2.58%  "
$ws.Range("D47").Value = "'0.634"
$ws.Range("E47").Value = "  -1.03%  "
$ws.Range("D48").Value = "'0.0253"
$ws.Range("E48").Value = "  -2.07%  "
$ws.Range("D49").Value = "'0.0991"
$ws.Range("E49").Value = "  -2.08%  "
$ws.Range("D50").Value = "'19.55"
$ws.Range("E50").Value = "  -3.05%  "
$ws.Range("D51").Value = "0.0₆0245"
$ws.Range("E51").Value = "  +6.09%  "
